$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the mixing-matrix values (rows 2-4, cols B-D) ---
# Row 2 ("0-20")
$ws.Range("B2").Value = 9.761902954
$ws.Range("C2").Value = 2.2706403989999999
$ws.Range("D2").Value = 0.29226176269999998

# Row 3 ("21-64")
$ws.Range("B3").Value = 6.6876464320000002
$ws.Range("C3").Value = 10.333415990000001
$ws.Range("D3").Value = 0.88528063099999998

# Row 4 ("65-100")
$ws.Range("B4").Value = 0.91568368420000001
$ws.Range("C4").Value = 1.01411684
$ws.Range("D4").Value = 1.2494830159999999

# These B/C/D cells move off the "italic/plain Times New Roman 12" style (s=1)
# back to the workbook default style (Normal / Calibri 12).
$ws.Range("B2:D4").Style = "Normal"

# A2 ("0-20" label) becomes bold, like the other header-ish label (matches
# the font used by A1/A6, i.e. bold Times New Roman 11).
$ws.Range("A2").Font.Name = "Times New Roman"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Bold = $true

# --- Update the selection / active cell recorded in the sheet view ---
$ws.Range("A15:E20").Select()
